$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$c3 = $cs.Colors(3)
$c3.RGB = 0x00FF00
Write-Output "done"
